# Append new Q&A rows (8-13) to the qa_database sheet, matching the
# target diff that extends the data range from A1:E7 to A1:E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(
        "What are the college's admission criteria?",
        "The college's admission criteria include academic performance, entrance exam scores, and extracurricular activities.",
        "कॉलेज के प्रवेश मानदंड क्या हैं?",
        "कॉलेज के प्रवेश मानदंडों में शैक्षणिक प्रदर्शन, प्रवेश परीक्षा के अंक, और सह-पाठ्यक्रम गतिविधियाँ शामिल हैं।",
        "admission criteria, college, academic performance, entrance exam"
    ),
    @(
        "What facilities are available on campus?",
        "The campus offers facilities such as a library, sports complex, canteen, and hostels.",
        "कैंपस में कौन-कौन सी सुविधाएँ उपलब्ध हैं?",
        "कैंपस में पुस्तकालय, खेल परिसर, कैफेटेरिया, और हॉस्टल जैसी सुविधाएँ उपलब्ध हैं।",
        "facilities, campus, library, sports complex, cafeteria, hostels"
    ),
    @(
        "What extracurricular activities does the college offer?",
        "The college offers various extracurricular activities including sports, cultural events, and technical workshops.",
        "कॉलेज कौन-कौन सी सह-पाठ्यक्रम गतिविधियाँ ऑफर करता है?",
        "कॉलेज विभिन्न सह-पाठ्यक्रम गतिविधियाँ ऑफर करता है जिनमें खेल, सांस्कृतिक कार्यक्रम, और तकनीकी कार्यशालाएँ शामिल हैं।",
        "extracurricular activities, college, sports, cultural events, technical clubs"
    ),
    @(
        "How can I contact the college administration?",
        "You can contact the college administration via phone, email, or by visiting the administrative office on campus.",
        "मैं कॉलेज प्रशासन से कैसे संपर्क कर सकता हूँ?",
        "आप फोन, ईमेल के माध्यम से या कैंपस में प्रशासनिक कार्यालय का दौरा करके कॉलेज प्रशासन से संपर्क कर सकते हैं।",
        "contact, college administration, phone, email, office"
    ),
    @(
        "Who is the addmission incharge of the college",
        "The admission incharge of the college is Mr. Atul Chakrawarti.",
        "कॉलेज के एडमिशन इंचार्ज कौन हैं?",
        "कॉलेज के एडमिशन इंचार्ज श्री अतुल चक्रवर्ती हैं।",
        "admission incharge, college, Mr. Atul Chakrawarti"
    ),
    @(
        "What is the location of the college?",
        "The college is located in Raipur, Chhattisgarh, India.",
        "कॉलेज का स्थान क्या है?",
        "कॉलेज रायपुर, छत्तीसगढ़, भारत में स्थित है।",
        "location, college, Raipur, Chhattisgarh"
    )
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
